$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 85.833336
$ws.Range("I8").Value = 67
$ws.Range("K8").Value = 201
$ws.Range("M8").Value = -62

# Row 15
$ws.Range("H15").Value = 688.451
$ws.Range("I15").Value = 688.451
$ws.Range("K15").Value = 2065.353
$ws.Range("M15").Value = -1896.353

# Row 58
$ws.Range("H58").Value = 1676.7778
$ws.Range("I58").Value = 886.375
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 2659.125
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -2509.125
$ws.Range("N58").Value = -24300

# Row 70
$ws.Range("H70").Value = 10461.529
$ws.Range("J70").Value = 14872.637
$ws.Range("L70").Value = 44617.911
$ws.Range("N70").Value = -45157.911

# Row 73
$ws.Range("H73").Value = 10461.529
$ws.Range("J73").Value = 14872.637
$ws.Range("L73").Value = 44617.911
$ws.Range("N73").Value = -46489.911

# Row 100
$ws.Range("H100").Value = 5137
$ws.Range("I100").Value = 1779.909
$ws.Range("K100").Value = 1779.909
$ws.Range("M100").Value = -1238.909

# Row 131
$ws.Range("H131").Value = 4936.6206
$ws.Range("I131").Value = 4150.619
$ws.Range("K131").Value = 12451.857
$ws.Range("M131").Value = -7411.857

# Row 132
$ws.Range("H132").Value = 993.7917
$ws.Range("I132").Value = 1013.54285
$ws.Range("K132").Value = 3040.62855
$ws.Range("M132").Value = -510.6285500000004

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3212.6667
$ws.Range("I2").Value = 2880.0715
$ws.Range("J2").Value = 3877.8572
$ws.Range("K2").Value = 2880.0715
$ws.Range("L2").Value = 3877.8572
$ws.Range("M2").Value = -2767.0715
$ws.Range("N2").Value = -4103.8572

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").ClearContents()

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("L82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("L85").ClearContents()

# Row 97
$ws.Range("H97").Value = 1845.7894
$ws.Range("I97").Value = 1841.9375
$ws.Range("J97").Value = 1866.3334
$ws.Range("K97").Value = 1841.9375
$ws.Range("L97").Value = 1866.3334
$ws.Range("M97").Value = -1345.9375
$ws.Range("N97").Value = -2858.3334

# Row 116
$ws.Range("H116").Value = 3212.6667
$ws.Range("I116").Value = 2880.0715
$ws.Range("J116").Value = 3877.8572
$ws.Range("K116").Value = 2880.0715
$ws.Range("L116").Value = 3877.8572
$ws.Range("M116").Value = -586.0715
$ws.Range("N116").Value = -8465.8572

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3212.6667
$ws.Range("I3").Value = 2880.0715
$ws.Range("J3").Value = 3877.8572
$ws.Range("K3").Value = 2880.0715
$ws.Range("L3").Value = 3877.8572
$ws.Range("M3").Value = -2766.0715
$ws.Range("N3").Value = -4105.8572

# Row 11
$ws.Range("H11").Value = 1931.75
$ws.Range("I11").Value = 237.25
$ws.Range("K11").Value = 237.25
$ws.Range("M11").Value = -97.25

# Row 86
$ws.Range("H86").Value = 2374
$ws.Range("I86").Value = 2500.4285
$ws.Range("J86").Value = 2152.75
$ws.Range("K86").Value = 2500.4285
$ws.Range("L86").Value = 2152.75
$ws.Range("M86").Value = -1377.4285
$ws.Range("N86").Value = -4398.75

# Row 89
$ws.Range("H89").Value = 2374
$ws.Range("I89").Value = 2500.4285
$ws.Range("J89").Value = 2152.75
$ws.Range("K89").Value = 12502.1425
$ws.Range("L89").Value = 10763.75
$ws.Range("M89").Value = -6886.1425
$ws.Range("N89").Value = -21995.75

# Row 94
$ws.Range("H94").Value = 1215.92
$ws.Range("I94").Value = 1189.3158
$ws.Range("K94").Value = 1189.3158
$ws.Range("M94").Value = -738.3158000000001

# Row 99
$ws.Range("H99").Value = 29018.15
$ws.Range("I99").Value = 37357.535
$ws.Range("K99").Value = 37357.535
$ws.Range("M99").Value = -35859.535

# Row 105
$ws.Range("H105").Value = 2665.8333
$ws.Range("I105").Value = 2665.8333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2665.8333
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = -918.8332999999998
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 7381.1816
$ws.Range("I132").Value = 8930.071
$ws.Range("J132").Value = 4670.625
$ws.Range("K132").Value = 26790.213
$ws.Range("L132").Value = 14011.875
$ws.Range("M132").Value = -24260.213
$ws.Range("N132").Value = -19071.875

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 250
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -659
$ws.Range("N46").Value = -15182

# Row 127
$ws.Range("H127").Value = 2000
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920

# Row 131
$ws.Range("H131").Value = 1697.4584
$ws.Range("I131").Value = 799.6667
$ws.Range("J131").Value = 1996.7222
$ws.Range("K131").Value = 2399.0001
$ws.Range("L131").Value = 5990.1666
$ws.Range("M131").Value = 2640.9999
$ws.Range("N131").Value = -16070.1666

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4215.8687
$ws.Range("I102").Value = 4325.8438
$ws.Range("J102").Value = 3629.3333
$ws.Range("K102").Value = 4325.8438
$ws.Range("L102").Value = 3629.3333
$ws.Range("M102").Value = -2703.8438
$ws.Range("N102").Value = -6873.3333

# Row 113
$ws.Range("H113").Value = 3550.4443
$ws.Range("I113").Value = 2101.6667
$ws.Range("K113").Value = 2101.6667
$ws.Range("M113").Value = 68.33329999999978

# Row 129
$ws.Range("H129").Value = 100000
$ws.Range("J129").Value = 100000
$ws.Range("L129").Value = 100000
$ws.Range("N129").Value = -110000

# Row 132
$ws.Range("H132").Value = 534556.9
$ws.Range("I132").Value = 1013869.25
$ws.Range("J132").Value = 11670.637
$ws.Range("K132").Value = 3041607.75
$ws.Range("L132").Value = 35011.911
$ws.Range("M132").Value = -3039077.75
$ws.Range("N132").Value = -40071.911

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3553.7778
$ws.Range("I68").Value = 2996.25
$ws.Range("J68").Value = 3999.8
$ws.Range("K68").Value = 2996.25
$ws.Range("L68").Value = 3999.8
$ws.Range("M68").Value = -2247.25
$ws.Range("N68").Value = -5497.8

# Row 71
$ws.Range("H71").Value = 3553.7778
$ws.Range("I71").Value = 2996.25
$ws.Range("J71").Value = 3999.8
$ws.Range("K71").Value = 14981.25
$ws.Range("L71").Value = 19999
$ws.Range("M71").Value = -11237.25
$ws.Range("N71").Value = -27487

# Row 93
$ws.Range("H93").Value = 2235
$ws.Range("I93").Value = 2316.4
$ws.Range("J93").Value = 2184.125
$ws.Range("K93").Value = 2316.4
$ws.Range("L93").Value = 2184.125
$ws.Range("M93").Value = -1068.4
$ws.Range("N93").Value = -4680.125

# Row 122
$ws.Range("H122").Value = 3299.2
$ws.Range("I122").Value = 3115.1072
$ws.Range("K122").Value = 9345.321599999999
$ws.Range("M122").Value = -6895.321599999999

# Row 132
$ws.Range("H132").Value = 1083394.5
$ws.Range("I132").Value = 1237450.9
$ws.Range("K132").Value = 3712352.7
$ws.Range("M132").Value = -3709822.7

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11144
